$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 06:22"

# Update India's stats (row 20)
$ws.Range("B20").Value = 21450
$ws.Range("C20").Value = 80
$ws.Range("D20").Value = 4373
$ws.Range("E20").Value = 16396

# Kazajistan's case counts overtook Hungria's, so the two countries swap rows (63/64)
# Row 63 becomes Kazajistan with its refreshed figures
$ws.Range("A63").Value = "Kazajistan"
$ws.Range("B63").Value = 2191
$ws.Range("C63").Value = 56
$ws.Range("D63").Value = 515
$ws.Range("E63").Value = 1657
$ws.Range("F63").Value = 29
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 19

# Row 64 becomes Hungria, keeping its previous (unchanged) figures
$ws.Range("A64").Value = "Hungria"
$ws.Range("B64").Value = 2168
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 295
$ws.Range("E64").Value = 1648
$ws.Range("F64").Value = 82
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 225
